$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The "Price" column (D) holds numeric-looking values (e.g. "1.00", "0.526")
# that must remain plain text, exactly like the source data. Temporarily force
# a Text number format on the whole column range before writing, then clear the
# formatting again afterwards so the cells keep their original (unstyled) look.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = '54.072.10'
$ws.Range("E2").Value = '  -3.53%  '
$ws.Range("D3").Value = '2.283.69'
$ws.Range("E3").Value = '  -3.59%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").Value = '492.61'
$ws.Range("E5").Value = '  -1.93%  '
$ws.Range("D6").Value = '128.17'
$ws.Range("E6").Value = '  -2.07%  '
$ws.Range("D7").Value = '1.00'
$ws.Range("E7").Value = '  +0.05%  '
$ws.Range("D8").Value = '0.526'
$ws.Range("E8").Value = '  -3.91%  '
$ws.Range("D9").Value = '2.288.67'
$ws.Range("E9").Value = '  -3.62%  '
$ws.Range("D10").Value = '0.0936'
$ws.Range("E10").Value = '  -3.68%  '
$ws.Range("E11").Value = '  -1.38%  '
$ws.Range("D12").Value = '4.76'
$ws.Range("E12").Value = '  +1.88%  '
$ws.Range("D13").Value = '0.317'
$ws.Range("E13").Value = '  -2.88%  '
$ws.Range("D14").Value = '2.693.45'
$ws.Range("E14").Value = '  -3.48%  '
$ws.Range("D15").Value = '21.34'
$ws.Range("E15").Value = '  -0.46%  '
$ws.Range("D16").Value = '54.096.79'
$ws.Range("E16").Value = '  -3.37%  '
$ws.Range("E17").Value = '  -1.84%  '
$ws.Range("D18").Value = '2.227.33'
$ws.Range("E18").Value = '  -2.90%  '
$ws.Range("D19").Value = '4.01'
$ws.Range("E19").Value = '  -0.04%  '
$ws.Range("D20").Value = '9.67'
$ws.Range("E20").Value = '  -3.42%  '
$ws.Range("D21").Value = '303.59'
$ws.Range("E21").Value = '  -1.12%  '
$ws.Range("D22").Value = '6.17'
$ws.Range("E22").Value = '  -1.22%  '
$ws.Range("E23").Value = '  -0.17%  '
$ws.Range("D24").Value = '63.91'
$ws.Range("E24").Value = '  -1.89%  '
$ws.Range("D25").Value = '0.999'
$ws.Range("E25").Value = '  +0.01%  '
$ws.Range("D26").Value = '0.367'
$ws.Range("E26").Value = '  -0.79%  '
$ws.Range("D27").Value = '0.143'
$ws.Range("E27").Value = '  -3.45%  '
$ws.Range("E28").Value = '  -1.60%  '
$ws.Range("D29").Value = '169.69'
$ws.Range("E29").Value = '  -1.81%  '
$ws.Range("E30").Value = '  -2.35%  '
$ws.Range("E31").Value = '  -1.31%  '
$ws.Range("D32").Value = '0.999'
$ws.Range("E32").Value = '  -0.09%  '
$ws.Range("E33").Value = '  +0.91%  '
$ws.Range("E34").Value = '  -0.10%  '
$ws.Range("D35").Value = '1.08'
$ws.Range("E35").Value = '  -1.82%  '
$ws.Range("D36").Value = '17.57'
$ws.Range("E36").Value = '  -0.20%  '
$ws.Range("E37").Value = '  -0.95%  '
$ws.Range("D38").Value = '0.850'
$ws.Range("E38").Value = '  +6.67%  '
$ws.Range("E39").Value = '  -4.09%  '
$ws.Range("D40").Value = '35.71'
$ws.Range("E40").Value = '  -0.82%  '
$ws.Range("E41").Value = '  -2.06%  '
$ws.Range("D42").Value = '0.368'
$ws.Range("E42").Value = '  -0.55%  '
$ws.Range("E43").Value = '  -0.75%  '
$ws.Range("B44").Value = 'Aave'
$ws.Range("C44").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D44").Value = '123.74'
$ws.Range("E44").Value = '  -5.59%  '
$ws.Range("B45").Value = 'RenderToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D45").Value = '4.70'
$ws.Range("E45").Value = '  -1.25%  '
$ws.Range("B46").Value = 'Stellar'
$ws.Range("C46").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D46").Value = '0.0883'
$ws.Range("E46").Value = '  -2.79%  '
$ws.Range("B47").Value = 'Mantle'
$ws.Range("C47").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D47").Value = '0.545'
$ws.Range("E47").Value = '  -3.58%  '
$ws.Range("D48").Value = '238.34'
$ws.Range("E48").Value = '  -2.22%  '
$ws.Range("E49").Value = '  -1.05%  '
$ws.Range("D50").Value = '0.0205'
$ws.Range("E50").Value = '  -1.21%  '
$ws.Range("D51").Value = '16.47'
$ws.Range("E51").Value = '  -2.57%  '

$priceRange.ClearFormats()
